$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> new DAMSLTag (col I) and DialogAct (col J) values,
# per the re-run of SGNN annotation after transcript clean-up.
$updates = @(
    @{ Row = 5;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 25; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 26; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 32; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 46; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 53; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 57; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 62; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 84; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 90; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.Tag
    $ws.Range("J" + $u.Row).Value = $u.Act
}
